$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '57.322.40'
$ws.Cells.Item(2, 5).Value = '  +0.50%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.430.77'
$ws.Cells.Item(3, 5).Value = '  -1.41%  '

$ws.Cells.Item(4, 5).Value = '  +0.34%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '489.95'
$ws.Cells.Item(5, 5).Value = '  -0.41%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '155.87'
$ws.Cells.Item(6, 5).Value = '  +1.60%  '

$ws.Cells.Item(7, 2).Value = 'USDC'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.996'
$ws.Cells.Item(7, 5).Value = '  -0.09%  '

$ws.Cells.Item(8, 2).Value = 'XRP'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.612'
$ws.Cells.Item(8, 5).Value = '  +19.47%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.462.01'
$ws.Cells.Item(9, 5).Value = '  +0.08%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.28'
$ws.Cells.Item(10, 5).Value = '  +10.51%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.101'
$ws.Cells.Item(11, 5).Value = '  -0.48%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.333'
$ws.Cells.Item(12, 5).Value = '  -0.71%  '

$ws.Cells.Item(13, 5).Value = '  +0.79%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.858.56'
$ws.Cells.Item(14, 5).Value = '  -0.97%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '57.023.38'
$ws.Cells.Item(15, 5).Value = '  -0.31%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '20.74'
$ws.Cells.Item(16, 5).Value = '  -1.43%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.0000134'
$ws.Cells.Item(17, 5).Value = '  -2.95%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.467.81'
$ws.Cells.Item(18, 5).Value = '  +0.07%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.66'
$ws.Cells.Item(19, 5).Value = '  +1.81%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '325.17'
$ws.Cells.Item(20, 5).Value = '  -0.40%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.04'
$ws.Cells.Item(21, 5).Value = '  -0.54%  '

$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.997'
$ws.Cells.Item(22, 5).Value = '  -0.15%  '

$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.92'
$ws.Cells.Item(23, 5).Value = '  +1.66%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '58.22'
$ws.Cells.Item(24, 5).Value = '  +0.37%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.402'
$ws.Cells.Item(25, 5).Value = '  -1.09%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.996'
$ws.Cells.Item(26, 5).Value = '  -0.27%  '

$ws.Cells.Item(27, 5).Value = '  -1.72%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.556.40'
$ws.Cells.Item(28, 5).Value = '  -0.15%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.31'
$ws.Cells.Item(29, 5).Value = '  -3.29%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0₃0802'
$ws.Cells.Item(30, 5).Value = '  -1.98%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.04%  '

$ws.Cells.Item(32, 2).Value = 'Monero'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '150.81'
$ws.Cells.Item(32, 5).Value = '  +0.41%  '

$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '18.75'
$ws.Cells.Item(33, 5).Value = '  +3.07%  '

$ws.Cells.Item(34, 5).Value = '  +0.01%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.33'
$ws.Cells.Item(35, 5).Value = '  +2.05%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.76'
$ws.Cells.Item(36, 5).Value = '  -0.34%  '

$ws.Cells.Item(37, 5).Value = '  -1.10%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.829'
$ws.Cells.Item(38, 5).Value = '  -7.96%  '

$ws.Cells.Item(39, 2).Value = 'OKB'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '34.28'
$ws.Cells.Item(39, 5).Value = '  +0.18%  '

$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.39'
$ws.Cells.Item(40, 5).Value = '  -1.00%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.55'
$ws.Cells.Item(41, 5).Value = '  +1.12%  '

$ws.Cells.Item(42, 5).Value = '  +4.92%  '

$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.996'
$ws.Cells.Item(43, 5).Value = '  +0.14%  '

$ws.Cells.Item(44, 2).Value = 'Bittensor'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '275.79'
$ws.Cells.Item(44, 5).Value = '  +3.35%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.603'
$ws.Cells.Item(45, 5).Value = '  -0.65%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0536'
$ws.Cells.Item(46, 5).Value = '  -4.22%  '

$ws.Cells.Item(47, 5).Value = '  +0.30%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0229'
$ws.Cells.Item(48, 5).Value = '  -0.51%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '4.64'
$ws.Cells.Item(49, 5).Value = '  -4.92%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '17.98'
$ws.Cells.Item(50, 5).Value = '  +0.87%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.886.51'
$ws.Cells.Item(51, 5).Value = '  +2.35%  '
